$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.975.97'
$ws.Range("E2").Value = '  +4.21%  '

$ws.Range("D3").Value = '2.421.00'
$ws.Range("E3").Value = '  +1.09%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.84'
$ws.Range("E5").Value = '  +3.78%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '102.66'
$ws.Range("E6").Value = '  +5.84%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.512'
$ws.Range("E7").Value = '  +1.54%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.524'
$ws.Range("E9").Value = '  +8.28%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.43'
$ws.Range("E10").Value = '  +1.55%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0799'
$ws.Range("E11").Value = '  +0.91%  '

$ws.Range("E12").Value = '  -2.18%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.21'
$ws.Range("E13").Value = '  -1.64%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.94'
$ws.Range("E14").Value = '  +1.35%  '

$ws.Range("D15").Value = '2.801.72'
$ws.Range("E15").Value = '  +1.67%  '

$ws.Range("D16").Value = '2.432.46'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.826'
$ws.Range("E17").Value = '  +1.39%  '

$ws.Range("D18").Value = '44.878.33'
$ws.Range("E18").Value = '  +3.90%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.16'
$ws.Range("E19").Value = '  +1.18%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.34'
$ws.Range("E20").Value = '  -0.43%  '

$ws.Range("E21").Value = '  +2.87%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.66'
$ws.Range("E22").Value = '  +0.58%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '242.41'
$ws.Range("E23").Value = '  +2.32%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.25'
$ws.Range("E24").Value = '  -0.06%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.49'
$ws.Range("E25").Value = '  +1.84%  '

$ws.Range("E26").Value = '  -0.09%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.16'
$ws.Range("E27").Value = '  +2.07%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.18'
$ws.Range("E28").Value = '  -8.05%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.48'
$ws.Range("E29").Value = '  +0.88%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '49.03'
$ws.Range("E30").Value = '  +1.98%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '32.71'
$ws.Range("E31").Value = '  +1.60%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.125'
$ws.Range("E32").Value = '  +9.55%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.86'
$ws.Range("E33").Value = '  +9.09%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.19'
$ws.Range("E34").Value = '  +1.68%  '

$ws.Range("E35").Value = '  +0.22%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0760'
$ws.Range("E36").Value = '  +2.42%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.87'
$ws.Range("E37").Value = '  +0.85%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.44'
$ws.Range("E38").Value = '  +2.20%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.85'
$ws.Range("E39").Value = '  -1.35%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '123.37'
$ws.Range("E40").Value = '  -5.70%  '

$ws.Range("B41").Value = 'WEMIXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.20'
$ws.Range("E41").Value = '  -3.00%  '

$ws.Range("B42").Value = 'Stellar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.109'
$ws.Range("E42").Value = '  +1.40%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.06'
$ws.Range("E43").Value = '  -0.06%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0288'
$ws.Range("E44").Value = '  +2.92%  '

$ws.Range("D45").Value = '1.933.39'
$ws.Range("E45").Value = '  -0.37%  '

$ws.Range("E46").Value = '  -2.29%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.90'
$ws.Range("E47").Value = '  +4.12%  '

$ws.Range("E48").Value = '  -0.24%  '

$ws.Range("E49").Value = '  +14.92%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '75.89'
$ws.Range("E50").Value = '  +5.58%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '53.09'
$ws.Range("E51").Value = '  +1.59%  '
